# edit.ps1 - apply the "edits to diagram and overview" change set to the
# architecture diagram deck.
#
# Coordinates in the OOXML are EMU (914400 EMU/inch, 12700 EMU/point).
# PowerPoint COM's Shape.Left/Top/Width/Height are in points (Single/float).
# This COM-interop shim truncates when converting points -> EMU, so a naive
# `emu / 12700.0` can land 1 EMU short. Adding half an EMU worth of points
# before the divide compensates for that truncation and reproduces the
# exact target EMU value.
function ToPt([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Big rounded-rectangle outline: widen it (cx 6650696 -> 6726926).
# ---------------------------------------------------------------------
$outline = $s.Shapes.Item(1)
$outline.Width = ToPt(6726926)

# ---------------------------------------------------------------------
# 2. "Graphic 22" (EKS icon, rId7/rId8 svg) moves up and to the right.
# ---------------------------------------------------------------------
$graphic22 = $s.Shapes.Item(9)
$graphic22.Left = ToPt(6283351)
$graphic22.Top = ToPt(2015)

# ---------------------------------------------------------------------
# 3. "pull metrics" textbox: shifts left & widens; run gets Consolas.
# ---------------------------------------------------------------------
$pullMetrics = $s.Shapes.Item(27)
$pullMetrics.Left = ToPt(2846978)
$pullMetrics.Width = ToPt(857927)
$pullMetrics.TextFrame.TextRange.Font.NameAscii = "Consolas"

# ---------------------------------------------------------------------
# 4. "push alerts" textbox: shifts left & widens; run gets Consolas.
# ---------------------------------------------------------------------
$pushAlerts = $s.Shapes.Item(35)
$pushAlerts.Left = ToPt(4812993)
$pushAlerts.Width = ToPt(801823)
$pushAlerts.TextFrame.TextRange.Font.NameAscii = "Consolas"

# ---------------------------------------------------------------------
# 5. New "Amazon Elastic Kubernetes Service (Amazon EKS)" label textbox,
#    appended to the slide. PowerPoint hands out shape ids from a
#    monotonically increasing, never-reused counter that skips any id
#    already present on the slide; the target id is 53, so burn through
#    the intervening ids with disposable textboxes first.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 30; $i++) {
    $scratch = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $scratch.Delete()
}

$labelLeft = ToPt(5633577)
$labelTop = ToPt(548361)
$labelWidth = ToPt(1852092)
$labelHeight = ToPt(369332)
$label = $s.Shapes.AddTextbox(1, $labelLeft, $labelTop, $labelWidth, $labelHeight)
$label.Name = "TextBox 9"

$label.TextFrame.WordWrap = -1
$label.TextFrame.AutoSize = 1

$label.Fill.Visible = 0
$label.Line.Visible = 0

$tr = $label.TextFrame.TextRange
$tr.Text = "Amazon Elastic Kubernetes Service (Amazon EKS)"
$tr.ParagraphFormat.Alignment = 2
$tr.ParagraphFormat.FarEastLineBreakControl = -1
$tr.ParagraphFormat.HangingPunctuation = -1

$lf = $tr.Font
$lf.Size = 9
$lf.NameAscii = "Consolas"
$lf.NameFarEast = "Amazon Ember"
$lf.NameComplexScript = "Arial"

Write-Output ("done; new shape id=" + $label.Id)
